# Generate Report for Handback
# Refresh the timestamp strings recorded on the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 0e0c301f... file.
$wsOverview.Range("G2").Value = "2016-08-15 11:01:56"

# zh-cn sheet: Correspond Handoff / Handback datetimes for the 0e0c301f... file.
$wsZhCn.Range("H2").Value = "2016-08-15 11:01:52"
$wsZhCn.Range("K2").Value = "2016-08-15 11:02:15"

# de-de sheet: "Latest HO Xliff Generate Date" (shares the same value as Overview!G2)
# and the Correspond Handback datetime for the 0e0c301f... file.
$wsDeDe.Range("H2").Value = "2016-08-15 11:01:56"
$wsDeDe.Range("K2").Value = "2016-08-15 11:02:22"
